$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @(44813, "KA01MR9978", "HONDA BRV",  "PMS",            "WORK DONE DELIVERED", 7415, "PAYTM"),
    @(44813, "KL40C4910",  "WAGON R",    "PMS",            "WORK DONE DELIVERED", 3835, "CARD"),
    @(44813, "KA19MB9767", "POLO",       "RUNNING REPAIR", "WORK DONE DELIVERED", 5500, "GPAY"),
    @(44813, "KA51N252",   "VERNA",      "HANDLE CHANGE",  "WORK DONE DELIVERED", 900,  "CASH"),
    @(44813, "KA03MK5127", "I10",        "HORN CHANGE",    "WORK DONE DELIVERED", 500,  "P PAY")
)

$startRow = 380
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$ws.Range("H380").Select()
